# Restored from revision of admin on 03/11/2021 07:54:51 AM.TEST Author: admin. Type: SAVE.
# Update cell C10 on the active sheet from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
